$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings
# (e.g. "1.00", "242.10") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '96.560.03'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '3.584.22'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '242.10'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '659.46'
$ws.Range("E6").Value = '  +1.45%  '
$ws.Range("D7").Value = '1.56'
$ws.Range("E7").Value = '  +7.21%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '1.05'
$ws.Range("E10").Value = '  +5.06%  '
$ws.Range("D11").Value = '3.582.44'
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = '43.55'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = '6.40'
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '4.252.06'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '96.501.81'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '0.0000258'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '3.586.50'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '7.75'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").Value = '12.63'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '17.89'
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("D22").Value = '0.493'
$ws.Range("E22").Value = '  +2.08%  '
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = '3.46'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '513.26'
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("D25").Value = '0.0000199'
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("D26").Value = '6.86'
$ws.Range("E26").Value = '  +3.39%  '
$ws.Range("D27").Value = '97.16'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '12.78'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").Value = '3.779.89'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").Value = '3.03'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("E31").Value = '  +7.64%  '
$ws.Range("D32").Value = '11.53'
$ws.Range("E32").Value = '  +2.60%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = '0.184'
$ws.Range("E34").Value = '  +4.89%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '31.77'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").Value = '0.568'
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("D38").Value = '8.49'
$ws.Range("E38").Value = '  +4.87%  '
$ws.Range("D39").Value = '596.32'
$ws.Range("E39").Value = '  +7.17%  '
$ws.Range("D40").Value = '1.60'
$ws.Range("E40").Value = '  +9.49%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.151'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '0.908'
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '1.85'
$ws.Range("E44").Value = '  +7.18%  '
$ws.Range("D45").Value = '5.75'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '2.28'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").Value = '34.24'
$ws.Range("E47").Value = '  +5.78%  '
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").Value = '23.55'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").Value = '3.65'
$ws.Range("E50").Value = '  +6.13%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '53.79'
$ws.Range("E51").Value = '  -1.05%  '

# Restore default (unstyled) cell style now that text values are locked in,
# matching the original workbook which had no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"
